$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp footer text (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 14:25"

# Swap the two country-name pairs whose rows carry identical statistics
# (Seychelles/Lesoto and Islas Malvinas/Groenlandia effectively traded places
# in the shared-string table, which shows up as the country label swapping
# between these adjacent rows).
$ws.Range("A184").Value = "Lesoto"
$ws.Range("A185").Value = "Seychelles"

$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# Refresh COVID-19 statistics for the updated countries

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3159671
$ws.Range("C4").Value = 739
$ws.Range("E4").Value = 1631527
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = 134888

# Row 19 - Alemania
$ws.Range("B19").Value = 198812
$ws.Range("C19").Value = 47
$ws.Range("E19").Value = 6097

# Row 34 - Emiratos Arabes Unidos
$ws.Range("B34").Value = 53577
$ws.Range("C34").Value = 532
$ws.Range("D34").Value = 43570
$ws.Range("E34").Value = 9679
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 328

# Row 36 - Kuwait
$ws.Range("B36").Value = 52840
$ws.Range("C36").Value = 833
$ws.Range("D36").Value = 42686
$ws.Range("E36").Value = 9772
$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 382

# Row 68 - Dinamarca
$ws.Range("B68").Value = 12916
$ws.Range("C68").Value = 16
$ws.Range("D68").Value = 12045
$ws.Range("E68").Value = 262

# Row 101 - Croacia
$ws.Range("B101").Value = 3416
$ws.Range("C101").Value = 91
$ws.Range("D101").Value = 2323
$ws.Range("E101").Value = 978
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 115

# Row 118 - Islandia
$ws.Range("B118").Value = 1882
$ws.Range("C118").Value = 2
$ws.Range("D118").Value = 1854
$ws.Range("E118").Value = 18

# Row 154 - Surinam
$ws.Range("B154").Value = 671
$ws.Range("C154").Value = 6
$ws.Range("D154").Value = 435
$ws.Range("E154").Value = 219

# Row 155 - Namibia
$ws.Range("B155").Value = 615
$ws.Range("C155").Value = 22
$ws.Range("E155").Value = 590
